$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.640398102298974
$ws.Range("C2").Value = 0.1587715242438037
$ws.Range("D2").Value = 0.07949264753537477
$ws.Range("E2").Value = 0.1316963412771841
$ws.Range("G2").Value = 0.3475388844066885
$ws.Range("H2").Value = 0.4842450376280283
$ws.Range("M2").Value = 0.3018572028508899
$ws.Range("N2").Value = 0.9455168968621663
$ws.Range("O2").Value = 1.598725780483932

$ws.Range("B3").Value = 0.5612502910020964
$ws.Range("C3").Value = 0.1426798511972152
$ws.Range("D3").Value = 0.07197596791050387
$ws.Range("E3").Value = 0.1250880501773466
$ws.Range("G3").Value = 0.3398822978595888
$ws.Range("H3").Value = 0.4847882341582306
$ws.Range("M3").Value = 0.2685275738221407
$ws.Range("N3").Value = 0.9555295310662544
$ws.Range("O3").Value = 1.583324172284108

$ws.Range("B4").Value = 0.512575376212709
$ws.Range("C4").Value = 0.1327393374131987
$ws.Range("D4").Value = 0.06739479089533518
$ws.Range("E4").Value = 0.1211461316568858
$ws.Range("G4").Value = 0.3355037645628443
$ws.Range("H4").Value = 0.4854148373558189
$ws.Range("M4").Value = 0.2481201624663427
$ws.Range("N4").Value = 0.9621064615574113
$ws.Range("O4").Value = 1.575154000857623

$ws.Range("B5").Value = 0.4927213306235103
$ws.Range("C5").Value = 0.1286735800856889
$ws.Range("D5").Value = 0.0655364917744663
$ws.Range("E5").Value = 0.1195685798303003
$ws.Range("G5").Value = 0.3338002977567953
$ws.Range("H5").Value = 0.4857438288041749
$ws.Range("M5").Value = 0.2398183476689226
$ws.Range("N5").Value = 0.9648945236705302
$ws.Range("O5").Value = 1.572147312993309

$ws.Range("B6").Value = 0.4894234847852772
$ws.Range("C6").Value = 0.1279975694385485
$ws.Range("D6").Value = 0.06522843993073479
$ws.Range("E6").Value = 0.1193083612383177
$ws.Range("G6").Value = 0.3335223107146845
$ws.Range("H6").Value = 0.4858029041029539
$ws.Range("M6").Value = 0.238440705955405
$ws.Range("N6").Value = 0.9653639968366292
$ws.Range("O6").Value = 1.571667523882667

$ws.Range("B7").Value = 0.51230769175109
$ws.Range("C7").Value = 0.1326845653080113
$ws.Range("D7").Value = 0.06736969454017583
$ws.Range("E7").Value = 0.1211247399210151
$ws.Range("G7").Value = 0.3354804641838314
$ws.Range("H7").Value = 0.4854189761308447
$ws.Range("M7").Value = 0.2480081430869134
$ws.Range("N7").Value = 0.962143625415159
$ws.Range("O7").Value = 1.575112145947088

$ws.Range("B8").Value = 0.6131246379628124
$ws.Range("C8").Value = 0.1532357235221298
$ws.Range("D8").Value = 0.07689381433226572
$ws.Range("E8").Value = 0.1293936533814559
$ws.Range("G8").Value = 0.3448317394008313
$ws.Range("H8").Value = 0.484371471107778
$ws.Range("M8").Value = 0.2903532870874486
$ws.Range("N8").Value = 0.948880186406349
$ws.Range("O8").Value = 1.593147823129868

$ws.Range("B9").Value = 0.8101765562271908
$ws.Range("C9").Value = 0.1930517069722839
$ws.Range("D9").Value = 0.09584236366829657
$ws.Range("E9").Value = 0.1465386960522252
$ws.Range("G9").Value = 0.3657450237356699
$ws.Range("H9").Value = 0.4846457114315257
$ws.Range("M9").Value = 0.3738491401714299
$ws.Range("N9").Value = 0.9262756263227772
$ws.Range("O9").Value = 1.638762624126542

$ws.Range("B10").Value = 0.9545253122140025
$ws.Range("C10").Value = 0.2220017199902884
$ws.Range("D10").Value = 0.1099327191611366
$ws.Range("E10").Value = 0.1597206078197857
$ws.Range("G10").Value = 0.3827035734410771
$ws.Range("H10").Value = 0.4862716434039243
$ws.Range("M10").Value = 0.4354840229545545
$ws.Range("N10").Value = 0.9117439608464579
$ws.Range("O10").Value = 1.678583339764344

$ws.Range("B11").Value = 1.020095658049172
$ws.Range("C11").Value = 0.2351047162476618
$ws.Range("D11").Value = 0.1163802107657688
$ws.Range("E11").Value = 0.1658486282549561
$ws.Range("G11").Value = 0.3907696477754286
$ws.Range("H11").Value = 0.4873217693036338
$ws.Range("M11").Value = 0.463589318609138
$ws.Range("N11").Value = 0.9055840192889022
$ws.Range("O11").Value = 1.698081875593488

$ws.Range("B12").Value = 1.044911068177726
$ws.Range("C12").Value = 0.2400567395843609
$ws.Range("D12").Value = 0.1188271591476706
$ws.Range("E12").Value = 0.168188363254032
$ws.Range("G12").Value = 0.3938749819930933
$ws.Range("H12").Value = 0.4877641486559838
$ws.Range("M12").Value = 0.4742418319334121
$ws.Range("N12").Value = 0.9033162343867716
$ws.Range("O12").Value = 1.705665419615201

$ws.Range("B13").Value = 1.039567294146252
$ws.Range("C13").Value = 0.2389906721785735
$ws.Range("D13").Value = 0.1182999238480278
$ws.Range("E13").Value = 0.1676836017123762
$ws.Range("G13").Value = 0.3932039241714875
$ws.Range("H13").Value = 0.4876668842780418
$ws.Range("M13").Value = 0.4719471939202435
$ws.Range("N13").Value = 0.9038017571975772
$ws.Range("O13").Value = 1.704023264546635

$ws.Range("B14").Value = 1.022137536520233
$ws.Range("C14").Value = 0.2355123198642843
$ws.Range("D14").Value = 0.1165814139972667
$ws.Range("E14").Value = 0.1660407335443139
$ws.Range("G14").Value = 0.3910241037244191
$ws.Range("H14").Value = 0.4873572674574831
$ws.Range("M14").Value = 0.4644655140853615
$ws.Range("N14").Value = 0.9053961471980188
$ws.Range("O14").Value = 1.698701767623561

$ws.Range("B15").Value = 1.011459366855945
$ws.Range("C15").Value = 0.2333804458842508
$ws.Range("D15").Value = 0.1155294833310307
$ws.Range("E15").Value = 0.1650369362502317
$ws.Range("G15").Value = 0.3896955376697093
$ws.Range("H15").Value = 0.4871734441683913
$ws.Range("M15").Value = 0.4598840248069536
$ws.Range("N15").Value = 0.9063812044474204
$ws.Range("O15").Value = 1.695468252073994

$ws.Range("B16").Value = 0.9502381419860626
$ws.Range("C16").Value = 0.2211440500902597
$ws.Range("D16").Value = 0.1095121193718569
$ws.Range("E16").Value = 0.1593227970805842
$ws.Range("G16").Value = 0.3821835391734254
$ws.Range("H16").Value = 0.4862092690884054
$ws.Range("M16").Value = 0.4336486273001725
$ws.Range("N16").Value = 0.9121555992712871
$ws.Range("O16").Value = 1.677336986424962

$ws.Range("B17").Value = 0.9126559195390769
$ws.Range("C17").Value = 0.2136202168173611
$ws.Range("D17").Value = 0.1058303217897816
$ws.Range("E17").Value = 0.1558512354229151
$ws.Range("G17").Value = 0.3776654493326106
$ws.Range("H17").Value = 0.4856973522291241
$ws.Range("M17").Value = 0.4175712730834817
$ws.Range("N17").Value = 0.9158134488413978
$ws.Range("O17").Value = 1.666569109622003

$ws.Range("B18").Value = 0.891030774987712
$ws.Range("C18").Value = 0.2092864605636748
$ws.Range("D18").Value = 0.1037161994924247
$ws.Range("E18").Value = 0.1538668460423622
$ws.Range("G18").Value = 0.3750998326613626
$ws.Range("H18").Value = 0.4854321307995235
$ws.Range("M18").Value = 0.4083303382897867
$ws.Range("N18").Value = 0.9179597547064446
$ws.Range("O18").Value = 1.660505918313532

$ws.Range("B19").Value = 0.8837073856364555
$ws.Range("C19").Value = 0.2078180594544676
$ws.Range("D19").Value = 0.1030010037821683
$ws.Range("E19").Value = 0.1531970804992682
$ws.Range("G19").Value = 0.3742368283621715
$ws.Range("H19").Value = 0.4853473476682524
$ws.Range("M19").Value = 0.4052026052634972
$ws.Range("N19").Value = 0.9186937388542518
$ws.Range("O19").Value = 1.65847536251141

$ws.Range("B20").Value = 0.9166575342306942
$ws.Range("C20").Value = 0.2144217898774059
$ws.Range("D20").Value = 0.1062218879888661
$ws.Range("E20").Value = 0.1562195077628772
$ws.Range("G20").Value = 0.3781429827489262
$ws.Range("H20").Value = 0.4857488220060446
$ws.Range("M20").Value = 0.4192820795261412
$ws.Range("N20").Value = 0.9154196750838395
$ws.Range("O20").Value = 1.667701886283339

$ws.Range("B21").Value = 1.027257484969311
$ws.Range("C21").Value = 0.2365342639947414
$ws.Range("D21").Value = 0.1170860348198346
$ws.Range("E21").Value = 0.1665227608193121
$ws.Range("G21").Value = 0.3916629861065104
$ws.Range("H21").Value = 0.487446995169222
$ws.Range("M21").Value = 0.4666628031821602
$ws.Range("N21").Value = 0.9049260756796471
$ws.Range("O21").Value = 1.700259388883097

$ws.Range("B22").Value = 1.099455030400122
$ws.Range("C22").Value = 0.2509288323749388
$ws.Range("D22").Value = 0.1242179929171243
$ws.Range("E22").Value = 0.1733684958091928
$ws.Range("G22").Value = 0.4007958526366622
$ws.Range("H22").Value = 0.4888175426294623
$ws.Range("M22").Value = 0.4976850680415481
$ws.Range("N22").Value = 0.8984459346727576
$ws.Range("O22").Value = 1.722703006344346

$ws.Range("B23").Value = 1.06093004899185
$ws.Range("C23").Value = 0.2432514919385085
$ws.Range("D23").Value = 0.120408640677141
$ws.Range("E23").Value = 0.1697044618955914
$ws.Range("G23").Value = 0.3958942055906363
$ws.Range("H23").Value = 0.4880621781655066
$ws.Range("M23").Value = 0.4811227657007748
$ws.Range("N23").Value = 0.9018698976576189
$ws.Range("O23").Value = 1.710617503587002

$ws.Range("B24").Value = 0.9148484621630928
$ws.Range("C24").Value = 0.2140594242367229
$ws.Range("D24").Value = 0.1060448528379396
$ws.Range("E24").Value = 0.1560529761715017
$ws.Range("G24").Value = 0.3779269905709413
$ws.Range("H24").Value = 0.485725461921021
$ws.Range("M24").Value = 0.4185086172382455
$ws.Range("N24").Value = 0.9155975651378299
$ws.Range("O24").Value = 1.667189361158819

$ws.Range("B25").Value = 0.7569413824425624
$ws.Range("C25").Value = 0.182333080663625
$ws.Range("D25").Value = 0.09068685325721049
$ws.Range("E25").Value = 0.1417990477271829
$ws.Range("G25").Value = 0.359809274189189
$ws.Range("H25").Value = 0.4843217362299583
$ws.Range("M25").Value = 0.3512111423917119
$ws.Range("N25").Value = 0.9320262352454165
$ws.Range("O25").Value = 1.625319328475058
